$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing "InferSent" task description (cell D19)
$ws.Range("D19").Value = "InferSent – Supervised Learning of Sentence Embeddings/Attention mechanism etc.: https://yashuseth.wordpress.com/2018/08/06/infersent-supervised-learning-of-sentence-embeddings/"

# Add new row 14: date 2023-08-14, new task about inferring pretrained sentence encoder
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 45152

$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = "infer pretrained sentence encoder without db, bc embedding is too big for maximum dense vector size"
$ws.Rows.Item(14).RowHeight = 34

# Row 15: add date 2023-08-15 and new task about huggingface sentence transformer,
# keeping the pre-existing D15 content untouched
$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 45153

$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "huggingface: init, save, load sentence transformer"

# Update selection to match new active cell
$ws.Range("B15").Select()
